{"js": "// Update Bible passage text to the v2 GetBible.net API wording.\n// The English KJV verse 2 drops a comma, and the Chinese Union\n// Traditional text gains inter-word spacing (word segmentation)\n// throughout Genesis 1 and Mark 1.\nconst pairs = [\n  [\" And the earth was without form, and void; and darkness was upon the face of the deep. And the Spirit of God moved upon the face of the waters.\", \" And the earth was without form and void; and darkness was upon the face of the deep. And the Spirit of God moved upon the face of the waters.\"],\n  [\" \\u8d77\\u521d\\uff0c\\u3000\\u795e\\u5275\\u9020\\u5929\\u5730\\u3002\", \" \\u8d77\\u521d\\uff0c\\u3000\\u795e \\u5275\\u9020 \\u5929 \\u5730\\u3002\"],\n  [\" \\u5730\\u662f\\u7a7a\\u865b\\u6df7\\u6c8c\\uff0c\\u6df5\\u9762\\u9ed1\\u6697\\uff1b\\u3000\\u795e\\u7684\\u9748\\u904b\\u884c\\u5728\\u6c34\\u9762\\u4e0a\\u3002\", \" \\u5730 \\u662f \\u7a7a\\u865b \\u6df7\\u6c8c\\uff0c\\u6df5\\u9762 \\u9ed1\\u6697\\uff1b\\u3000\\u795e \\u7684\\u9748 \\u904b\\u884c \\u5728 \\u6c34 \\u9762\\u4e0a\\u3002\"],\n  [\" \\u795e\\u8aaa\\uff1a\\u300c\\u8981\\u6709\\u5149\\u300d\\uff0c\\u5c31\\u6709\\u4e86\\u5149\\u3002\", \" \\u795e \\u8aaa\\uff1a\\u300c\\u8981\\u6709 \\u5149\\u300d\\uff0c\\u5c31\\u6709\\u4e86 \\u5149\\u3002\"],\n  [\" \\u795e\\u770b\\u5149\\u662f\\u597d\\u7684\\uff0c\\u5c31\\u628a\\u5149\\u6697\\u5206\\u958b\\u4e86\\u3002\", \" \\u795e \\u770b \\u5149 \\u662f\\u597d\\u7684\\uff0c \\u5c31\\u628a \\u5149 \\u6697 \\u5206\\u958b\\u4e86\\u3002\"],\n  [\" \\u795e\\u5c31\\u7167\\u8457\\u81ea\\u5df1\\u7684\\u5f62\\u50cf\\u9020\\u4eba\\uff0c\\u4e43\\u662f\\u7167\\u8457\\u4ed6\\u7684\\u5f62\\u50cf\\u9020\\u7537\\u9020\\u5973\\u3002\", \" \\u795e \\u5c31\\u7167\\u8457\\u81ea\\u5df1\\u7684\\u5f62\\u50cf \\u9020 \\u4eba\\uff0c\\u4e43\\u662f\\u7167\\u8457\\u4ed6 \\u7684\\u5f62\\u50cf \\u9020 \\u7537 \\u9020 \\u5973\\u3002\"],\n  [\" \\u795e\\u5c31\\u8cdc\\u798f\\u7d66\\u4ed6\\u5011\\uff0c\\u53c8\\u5c0d\\u4ed6\\u5011\\u8aaa\\uff1a\\u300c\\u8981\\u751f\\u990a\\u773e\\u591a\\uff0c\\u904d\\u6eff\\u5730\\u9762\\uff0c\\u6cbb\\u7406\\u9019\\u5730\\uff0c\\u4e5f\\u8981\\u7ba1\\u7406\\u6d77\\u88e1\\u7684\\u9b5a\\u3001\\u7a7a\\u4e2d\\u7684\\u9ce5\\uff0c\\u548c\\u5730\\u4e0a\\u5404\\u6a23\\u884c\\u52d5\\u7684\\u6d3b\\u7269\\u3002\\u300d\", \" \\u795e \\u5c31\\u8cdc\\u798f \\u7d66\\u4ed6\\u5011\\uff0c \\u53c8\\u5c0d\\u4ed6\\u5011 \\u8aaa\\uff1a\\u300c\\u8981\\u751f\\u990a \\u773e\\u591a\\uff0c\\u904d\\u6eff \\u5730\\u9762\\uff0c\\u6cbb\\u7406\\u9019\\u5730\\uff0c\\u4e5f\\u8981\\u7ba1\\u7406 \\u6d77 \\u88e1\\u7684\\u9b5a\\u3001\\u7a7a \\u4e2d\\u7684\\u9ce5\\uff0c\\u548c\\u5730\\u4e0a \\u5404\\u6a23 \\u884c\\u52d5\\u7684 \\u6d3b\\u7269\\u3002\\u300d\"],\n  [\" \\u795e\\u8aaa\\uff1a\\u300c\\u770b\\u54ea\\uff0c\\u6211\\u5c07\\u904d\\u5730\\u4e0a\\u4e00\\u5207\\u7d50\\u7a2e\\u5b50\\u7684\\u83dc\\u852c\\u548c\\u4e00\\u5207\\u6a39\\u4e0a\\u6240\\u7d50\\u6709\\u6838\\u7684\\u679c\\u5b50\\u5168\\u8cdc\\u7d66\\u4f60\\u5011\\u4f5c\\u98df\\u7269\\u3002\", \" \\u795e \\u8aaa\\uff1a\\u300c\\u770b\\u54ea\\uff0c\\u6211\\u5c07 \\u904d \\u5730 \\u4e0a \\u4e00\\u5207 \\u7d50 \\u7a2e\\u5b50 \\u7684\\u83dc\\u852c \\u548c \\u4e00\\u5207 \\u6a39 \\u4e0a \\u6240 \\u7d50 \\u6709\\u6838 \\u7684\\u679c\\u5b50 \\u5168\\u8cdc \\u7d66\\u4f60\\u5011 \\u4f5c \\u98df\\u7269\\u3002\"],\n  [\" \\u81f3\\u65bc\\u5730\\u4e0a\\u7684\\u8d70\\u7378\\u548c\\u7a7a\\u4e2d\\u7684\\u98db\\u9ce5\\uff0c\\u4e26\\u5404\\u6a23\\u722c\\u5728\\u5730\\u4e0a\\u6709\\u751f\\u547d\\u7684\\u7269\\uff0c\\u6211\\u5c07\\u9752\\u8349\\u8cdc\\u7d66\\u7260\\u5011\\u4f5c\\u98df\\u7269\\u3002\\u300d\\u4e8b\\u5c31\\u9019\\u6a23\\u6210\\u4e86\\u3002\", \" \\u81f3\\u65bc \\u5730\\u4e0a \\u7684\\u8d70\\u7378\\u548c \\u7a7a\\u4e2d \\u7684\\u98db\\u9ce5\\uff0c\\u4e26\\u5404\\u6a23 \\u722c \\u5728 \\u5730\\u4e0a\\u6709 \\u751f\\u547d \\u7684\\u7269\\uff0c\\u6211\\u5c07 \\u9752 \\u8349 \\u8cdc\\u7d66\\u7260\\u5011\\u4f5c\\u98df\\u7269\\u3002\\u300d\\u4e8b\\u5c31\\u9019\\u6a23 \\u6210\\u4e86\\u3002\"],\n  [\" \\u795e\\u770b\\u8457\\u4e00\\u5207\\u6240\\u9020\\u7684\\u90fd\\u751a\\u597d\\u3002\\u6709\\u665a\\u4e0a\\uff0c\\u6709\\u65e9\\u6668\\uff0c\\u662f\\u7b2c\\u516d\\u65e5\\u3002\", \" \\u795e \\u770b\\u8457 \\u4e00\\u5207 \\u6240 \\u9020\\u7684 \\u90fd\\u751a \\u597d\\u3002\\u6709 \\u665a\\u4e0a\\uff0c\\u6709 \\u65e9\\u6668\\uff0c\\u662f\\u7b2c\\u516d \\u65e5\\u3002\"],\n  [\" \\u795e\\u7684\\u5152\\u5b50\\uff0c\\u8036\\u7a4c\\u57fa\\u7763\\u798f\\u97f3\\u7684\\u8d77\\u982d\\u3002\", \" \\u795e\\u7684 \\u5152\\u5b50\\uff0c\\u8036\\u7a4c \\u57fa\\u7763 \\u798f\\u97f3\\u7684 \\u8d77\\u982d\\u3002\"],\n  [\" \\u6b63\\u5982\\u5148\\u77e5\\u4ee5\\u8cfd\\u4e9e\\uff08\\u6709\\u53e4\\u5377\\u6c92\\u6709\\u4ee5\\u8cfd\\u4e9e\\u4e09\\u500b\\u5b57\\uff09\\u66f8\\u4e0a\\u8a18\\u8457\\u8aaa\\uff1a\\u770b\\u54ea\\uff0c\\u6211\\u8981\\u5dee\\u9063\\u6211\\u7684\\u4f7f\\u8005\\u5728\\u4f60\\u524d\\u9762\\uff0c\\u9810\\u5099\\u9053\\u8def\\u3002\", \" \\u6b63\\u5982 \\u5148\\u77e5 \\u4ee5\\u8cfd\\u4e9e\\uff08\\u6709\\u53e4\\u5377\\u6c92\\u6709\\u4ee5\\u8cfd\\u4e9e\\u4e09\\u500b\\u5b57\\uff09\\u66f8\\u4e0a\\u8a18\\u8457\\u8aaa\\uff1a\\u770b\\u54ea\\uff0c\\u6211\\u8981\\u5dee\\u9063 \\u6211\\u7684 \\u4f7f\\u8005 \\u5728 \\u4f60 \\u524d\\u9762\\uff0c \\u9810\\u5099 \\u9053\\u8def\\u3002\"],\n  [\" \\u5728\\u66e0\\u91ce\\u6709\\u4eba\\u8072\\u558a\\u8457\\u8aaa\\uff1a\\u9810\\u5099\\u4e3b\\u7684\\u9053\\uff0c\\u4fee\\u76f4\\u4ed6\\u7684\\u8def\\u3002\", \" \\u5728 \\u66e0\\u91ce \\u6709\\u4eba\\u8072 \\u558a\\u8457\\u8aaa\\uff1a\\u9810\\u5099 \\u4e3b\\u7684 \\u9053\\uff0c\\u4fee \\u76f4 \\u4ed6\\u7684 \\u8def\\u3002\"],\n  [\" \\u7d04\\u7ff0\\u4e0b\\u76e3\\u4ee5\\u5f8c\\uff0c\\u8036\\u7a4c\\u4f86\\u5230\\u52a0\\u5229\\u5229\\uff0c\\u5ba3\\u50b3\\u3000\\u795e\\u7684\\u798f\\u97f3\\uff0c\", \" \\u7d04\\u7ff0 \\u4e0b\\u76e3 \\u4ee5\\u5f8c\\uff0c\\u8036\\u7a4c \\u4f86 \\u5230 \\u52a0\\u5229\\u5229\\uff0c\\u5ba3\\u50b3\\u3000\\u795e\\u7684 \\u798f\\u97f3\\uff0c\"],\n  [\" \\u8aaa\\uff1a\\u300c\\u65e5\\u671f\\u6eff\\u4e86\\uff0c\\u3000\\u795e\\u7684\\u570b\\u8fd1\\u4e86\\u3002\\u4f60\\u5011\\u7576\\u6094\\u6539\\uff0c\\u4fe1\\u798f\\u97f3\\uff01\\u300d\", \" \\u8aaa\\uff1a\\u300c \\u65e5\\u671f \\u6eff\\u4e86\\uff0c\\u3000\\u795e\\u7684 \\u570b \\u8fd1\\u4e86\\u3002\\u4f60\\u5011\\u7576\\u6094\\u6539\\uff0c \\u4fe1 \\u798f\\u97f3\\uff01\\u300d\"],\n  [\" \\u8036\\u7a4c\\u9806\\u8457\\u52a0\\u5229\\u5229\\u7684\\u6d77\\u908a\\u8d70\\uff0c\\u770b\\u898b\\u897f\\u9580\\u548c\\u897f\\u9580\\u7684\\u5144\\u5f1f\\u5b89\\u5f97\\u70c8\\u5728\\u6d77\\u88e1\\u6492\\u7db2\\uff1b\\u4ed6\\u5011\\u672c\\u662f\\u6253\\u9b5a\\u7684\\u3002\", \" \\u8036\\u7a4c\\u9806\\u8457 \\u52a0\\u5229\\u5229\\u7684 \\u6d77\\u908a \\u8d70\\uff0c\\u770b\\u898b \\u897f\\u9580 \\u548c \\u897f\\u9580\\u7684 \\u5144\\u5f1f \\u5b89\\u5f97\\u70c8 \\u5728 \\u6d77 \\u88e1\\u6492\\u7db2\\uff1b\\u4ed6\\u5011\\u672c \\u662f \\u6253\\u9b5a\\u7684\\u3002\"],\n  [\" \\u8036\\u7a4c\\u5c0d\\u4ed6\\u5011\\u8aaa\\uff1a\\u300c\\u4f86\\u8ddf\\u5f9e\\u6211\\uff0c\\u6211\\u8981\\u53eb\\u4f60\\u5011\\u5f97\\u4eba\\u5982\\u5f97\\u9b5a\\u4e00\\u6a23\\u3002\\u300d\", \" \\u8036\\u7a4c \\u5c0d\\u4ed6\\u5011 \\u8aaa\\uff1a\\u300c\\u4f86 \\u8ddf\\u5f9e \\u6211\\uff0c \\u6211\\u8981\\u53eb \\u4f60\\u5011 \\u5f97\\u4eba \\u5982\\u5f97\\u9b5a\\u4e00\\u6a23\\u3002\\u300d\"],\n  [\" \\u4ed6\\u5011\\u5c31\\u7acb\\u523b\\u6368\\u4e86\\u7db2\\uff0c\\u8ddf\\u5f9e\\u4e86\\u4ed6\\u3002\", \" \\u4ed6\\u5011\\u5c31 \\u7acb\\u523b \\u6368\\u4e86 \\u7db2\\uff0c\\u8ddf\\u5f9e\\u4e86 \\u4ed6\\u3002\"],\n  [\" \\u8036\\u7a4c\\u7a0d\\u5f80\\u524d\\u8d70\\uff0c\\u53c8\\u898b\\u897f\\u5e87\\u592a\\u7684\\u5152\\u5b50\\u96c5\\u5404\\u548c\\u96c5\\u5404\\u7684\\u5144\\u5f1f\\u7d04\\u7ff0\\u5728\\u8239\\u4e0a\\u88dc\\u7db2\\u3002\", \" \\u8036\\u7a4c\\u7a0d \\u5f80\\u524d\\u8d70\\uff0c\\u53c8\\u898b \\u897f\\u5e87\\u592a\\u7684 \\u5152\\u5b50\\u96c5\\u5404 \\u548c \\u96c5\\u5404\\u7684 \\u5144\\u5f1f \\u7d04\\u7ff0 \\u5728 \\u8239\\u4e0a \\u88dc \\u7db2\\u3002\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of pairs) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(\"Text not found: \" + oldText);\n  }\n\n  for (const range of results.items) {\n    range.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Update Bible passage text to the v2 GetBible.net API wording.\n# The English KJV verse 2 drops a comma, and the Chinese Union\n# Traditional text gains inter-word spacing (word segmentation)\n# throughout Genesis 1 and Mark 1.\n$pairs = @(\n  @{old=\" And the earth was without form, and void; and darkness was upon the face of the deep. And the Spirit of God moved upon the face of the waters.\"; new=\" And the earth was without form and void; and darkness was upon the face of the deep. And the Spirit of God moved upon the face of the waters.\"},\n  @{old=\" \u8d77\u521d\uff0c\u3000\u795e\u5275\u9020\u5929\u5730\u3002\"; new=\" \u8d77\u521d\uff0c\u3000\u795e \u5275\u9020 \u5929 \u5730\u3002\"},\n  @{old=\" \u5730\u662f\u7a7a\u865b\u6df7\u6c8c\uff0c\u6df5\u9762\u9ed1\u6697\uff1b\u3000\u795e\u7684\u9748\u904b\u884c\u5728\u6c34\u9762\u4e0a\u3002\"; new=\" \u5730 \u662f \u7a7a\u865b \u6df7\u6c8c\uff0c\u6df5\u9762 \u9ed1\u6697\uff1b\u3000\u795e \u7684\u9748 \u904b\u884c \u5728 \u6c34 \u9762\u4e0a\u3002\"},\n  @{old=\" \u795e\u8aaa\uff1a\u300c\u8981\u6709\u5149\u300d\uff0c\u5c31\u6709\u4e86\u5149\u3002\"; new=\" \u795e \u8aaa\uff1a\u300c\u8981\u6709 \u5149\u300d\uff0c\u5c31\u6709\u4e86 \u5149\u3002\"},\n  @{old=\" \u795e\u770b\u5149\u662f\u597d\u7684\uff0c\u5c31\u628a\u5149\u6697\u5206\u958b\u4e86\u3002\"; new=\" \u795e \u770b \u5149 \u662f\u597d\u7684\uff0c \u5c31\u628a \u5149 \u6697 \u5206\u958b\u4e86\u3002\"},\n  @{old=\" \u795e\u5c31\u7167\u8457\u81ea\u5df1\u7684\u5f62\u50cf\u9020\u4eba\uff0c\u4e43\u662f\u7167\u8457\u4ed6\u7684\u5f62\u50cf\u9020\u7537\u9020\u5973\u3002\"; new=\" \u795e \u5c31\u7167\u8457\u81ea\u5df1\u7684\u5f62\u50cf \u9020 \u4eba\uff0c\u4e43\u662f\u7167\u8457\u4ed6 \u7684\u5f62\u50cf \u9020 \u7537 \u9020 \u5973\u3002\"},\n  @{old=\" \u795e\u5c31\u8cdc\u798f\u7d66\u4ed6\u5011\uff0c\u53c8\u5c0d\u4ed6\u5011\u8aaa\uff1a\u300c\u8981\u751f\u990a\u773e\u591a\uff0c\u904d\u6eff\u5730\u9762\uff0c\u6cbb\u7406\u9019\u5730\uff0c\u4e5f\u8981\u7ba1\u7406\u6d77\u88e1\u7684\u9b5a\u3001\u7a7a\u4e2d\u7684\u9ce5\uff0c\u548c\u5730\u4e0a\u5404\u6a23\u884c\u52d5\u7684\u6d3b\u7269\u3002\u300d\"; new=\" \u795e \u5c31\u8cdc\u798f \u7d66\u4ed6\u5011\uff0c \u53c8\u5c0d\u4ed6\u5011 \u8aaa\uff1a\u300c\u8981\u751f\u990a \u773e\u591a\uff0c\u904d\u6eff \u5730\u9762\uff0c\u6cbb\u7406\u9019\u5730\uff0c\u4e5f\u8981\u7ba1\u7406 \u6d77 \u88e1\u7684\u9b5a\u3001\u7a7a \u4e2d\u7684\u9ce5\uff0c\u548c\u5730\u4e0a \u5404\u6a23 \u884c\u52d5\u7684 \u6d3b\u7269\u3002\u300d\"},\n  @{old=\" \u795e\u8aaa\uff1a\u300c\u770b\u54ea\uff0c\u6211\u5c07\u904d\u5730\u4e0a\u4e00\u5207\u7d50\u7a2e\u5b50\u7684\u83dc\u852c\u548c\u4e00\u5207\u6a39\u4e0a\u6240\u7d50\u6709\u6838\u7684\u679c\u5b50\u5168\u8cdc\u7d66\u4f60\u5011\u4f5c\u98df\u7269\u3002\"; new=\" \u795e \u8aaa\uff1a\u300c\u770b\u54ea\uff0c\u6211\u5c07 \u904d \u5730 \u4e0a \u4e00\u5207 \u7d50 \u7a2e\u5b50 \u7684\u83dc\u852c \u548c \u4e00\u5207 \u6a39 \u4e0a \u6240 \u7d50 \u6709\u6838 \u7684\u679c\u5b50 \u5168\u8cdc \u7d66\u4f60\u5011 \u4f5c \u98df\u7269\u3002\"},\n  @{old=\" \u81f3\u65bc\u5730\u4e0a\u7684\u8d70\u7378\u548c\u7a7a\u4e2d\u7684\u98db\u9ce5\uff0c\u4e26\u5404\u6a23\u722c\u5728\u5730\u4e0a\u6709\u751f\u547d\u7684\u7269\uff0c\u6211\u5c07\u9752\u8349\u8cdc\u7d66\u7260\u5011\u4f5c\u98df\u7269\u3002\u300d\u4e8b\u5c31\u9019\u6a23\u6210\u4e86\u3002\"; new=\" \u81f3\u65bc \u5730\u4e0a \u7684\u8d70\u7378\u548c \u7a7a\u4e2d \u7684\u98db\u9ce5\uff0c\u4e26\u5404\u6a23 \u722c \u5728 \u5730\u4e0a\u6709 \u751f\u547d \u7684\u7269\uff0c\u6211\u5c07 \u9752 \u8349 \u8cdc\u7d66\u7260\u5011\u4f5c\u98df\u7269\u3002\u300d\u4e8b\u5c31\u9019\u6a23 \u6210\u4e86\u3002\"},\n  @{old=\" \u795e\u770b\u8457\u4e00\u5207\u6240\u9020\u7684\u90fd\u751a\u597d\u3002\u6709\u665a\u4e0a\uff0c\u6709\u65e9\u6668\uff0c\u662f\u7b2c\u516d\u65e5\u3002\"; new=\" \u795e \u770b\u8457 \u4e00\u5207 \u6240 \u9020\u7684 \u90fd\u751a \u597d\u3002\u6709 \u665a\u4e0a\uff0c\u6709 \u65e9\u6668\uff0c\u662f\u7b2c\u516d \u65e5\u3002\"},\n  @{old=\" \u795e\u7684\u5152\u5b50\uff0c\u8036\u7a4c\u57fa\u7763\u798f\u97f3\u7684\u8d77\u982d\u3002\"; new=\" \u795e\u7684 \u5152\u5b50\uff0c\u8036\u7a4c \u57fa\u7763 \u798f\u97f3\u7684 \u8d77\u982d\u3002\"},\n  @{old=\" \u6b63\u5982\u5148\u77e5\u4ee5\u8cfd\u4e9e\uff08\u6709\u53e4\u5377\u6c92\u6709\u4ee5\u8cfd\u4e9e\u4e09\u500b\u5b57\uff09\u66f8\u4e0a\u8a18\u8457\u8aaa\uff1a\u770b\u54ea\uff0c\u6211\u8981\u5dee\u9063\u6211\u7684\u4f7f\u8005\u5728\u4f60\u524d\u9762\uff0c\u9810\u5099\u9053\u8def\u3002\"; new=\" \u6b63\u5982 \u5148\u77e5 \u4ee5\u8cfd\u4e9e\uff08\u6709\u53e4\u5377\u6c92\u6709\u4ee5\u8cfd\u4e9e\u4e09\u500b\u5b57\uff09\u66f8\u4e0a\u8a18\u8457\u8aaa\uff1a\u770b\u54ea\uff0c\u6211\u8981\u5dee\u9063 \u6211\u7684 \u4f7f\u8005 \u5728 \u4f60 \u524d\u9762\uff0c \u9810\u5099 \u9053\u8def\u3002\"},\n  @{old=\" \u5728\u66e0\u91ce\u6709\u4eba\u8072\u558a\u8457\u8aaa\uff1a\u9810\u5099\u4e3b\u7684\u9053\uff0c\u4fee\u76f4\u4ed6\u7684\u8def\u3002\"; new=\" \u5728 \u66e0\u91ce \u6709\u4eba\u8072 \u558a\u8457\u8aaa\uff1a\u9810\u5099 \u4e3b\u7684 \u9053\uff0c\u4fee \u76f4 \u4ed6\u7684 \u8def\u3002\"},\n  @{old=\" \u7d04\u7ff0\u4e0b\u76e3\u4ee5\u5f8c\uff0c\u8036\u7a4c\u4f86\u5230\u52a0\u5229\u5229\uff0c\u5ba3\u50b3\u3000\u795e\u7684\u798f\u97f3\uff0c\"; new=\" \u7d04\u7ff0 \u4e0b\u76e3 \u4ee5\u5f8c\uff0c\u8036\u7a4c \u4f86 \u5230 \u52a0\u5229\u5229\uff0c\u5ba3\u50b3\u3000\u795e\u7684 \u798f\u97f3\uff0c\"},\n  @{old=\" \u8aaa\uff1a\u300c\u65e5\u671f\u6eff\u4e86\uff0c\u3000\u795e\u7684\u570b\u8fd1\u4e86\u3002\u4f60\u5011\u7576\u6094\u6539\uff0c\u4fe1\u798f\u97f3\uff01\u300d\"; new=\" \u8aaa\uff1a\u300c \u65e5\u671f \u6eff\u4e86\uff0c\u3000\u795e\u7684 \u570b \u8fd1\u4e86\u3002\u4f60\u5011\u7576\u6094\u6539\uff0c \u4fe1 \u798f\u97f3\uff01\u300d\"},\n  @{old=\" \u8036\u7a4c\u9806\u8457\u52a0\u5229\u5229\u7684\u6d77\u908a\u8d70\uff0c\u770b\u898b\u897f\u9580\u548c\u897f\u9580\u7684\u5144\u5f1f\u5b89\u5f97\u70c8\u5728\u6d77\u88e1\u6492\u7db2\uff1b\u4ed6\u5011\u672c\u662f\u6253\u9b5a\u7684\u3002\"; new=\" \u8036\u7a4c\u9806\u8457 \u52a0\u5229\u5229\u7684 \u6d77\u908a \u8d70\uff0c\u770b\u898b \u897f\u9580 \u548c \u897f\u9580\u7684 \u5144\u5f1f \u5b89\u5f97\u70c8 \u5728 \u6d77 \u88e1\u6492\u7db2\uff1b\u4ed6\u5011\u672c \u662f \u6253\u9b5a\u7684\u3002\"},\n  @{old=\" \u8036\u7a4c\u5c0d\u4ed6\u5011\u8aaa\uff1a\u300c\u4f86\u8ddf\u5f9e\u6211\uff0c\u6211\u8981\u53eb\u4f60\u5011\u5f97\u4eba\u5982\u5f97\u9b5a\u4e00\u6a23\u3002\u300d\"; new=\" \u8036\u7a4c \u5c0d\u4ed6\u5011 \u8aaa\uff1a\u300c\u4f86 \u8ddf\u5f9e \u6211\uff0c \u6211\u8981\u53eb \u4f60\u5011 \u5f97\u4eba \u5982\u5f97\u9b5a\u4e00\u6a23\u3002\u300d\"},\n  @{old=\" \u4ed6\u5011\u5c31\u7acb\u523b\u6368\u4e86\u7db2\uff0c\u8ddf\u5f9e\u4e86\u4ed6\u3002\"; new=\" \u4ed6\u5011\u5c31 \u7acb\u523b \u6368\u4e86 \u7db2\uff0c\u8ddf\u5f9e\u4e86 \u4ed6\u3002\"},\n  @{old=\" \u8036\u7a4c\u7a0d\u5f80\u524d\u8d70\uff0c\u53c8\u898b\u897f\u5e87\u592a\u7684\u5152\u5b50\u96c5\u5404\u548c\u96c5\u5404\u7684\u5144\u5f1f\u7d04\u7ff0\u5728\u8239\u4e0a\u88dc\u7db2\u3002\"; new=\" \u8036\u7a4c\u7a0d \u5f80\u524d\u8d70\uff0c\u53c8\u898b \u897f\u5e87\u592a\u7684 \u5152\u5b50\u96c5\u5404 \u548c \u96c5\u5404\u7684 \u5144\u5f1f \u7d04\u7ff0 \u5728 \u8239\u4e0a \u88dc \u7db2\u3002\"},\n)\n\n$d = $word.ActiveDocument\n\nforeach ($p in $pairs) {\n  $rng = $d.Content\n  $rng.Find.ClearFormatting()\n  $rng.Find.Replacement.ClearFormatting()\n  $found = $rng.Find.Execute($p.old, $false, $false, $false, $false, $false, $true, 1, $false, $p.new, 2)\n  if (-not $found) {\n    throw \"Text not found: $($p.old)\"\n  }\n}\n"}
